$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "Usb-tikulla tai itse valitsemallasi turvallisella tavalla." gets
# explicit fi-FI language formatting (no visible text change).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Usb-tikulla", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
if ($found) {
    $r.LanguageID = "fi-FI"
}

$r = $d.Content
$found = $r.Find.Execute(" tai itse valitsemallasi turvallisella tavalla.", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    " tai itse valitsemallasi turvallisella tavalla.", 2)
if ($found) {
    $r.LanguageID = "fi-FI"
}

# ---------------------------------------------------------------------------
# Hunk 2: "...lla tavalla." -> "...lla tavalla (Signal jne.)." in the
# "Lisäksi sovellus..." paragraph.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(3).Range
$found = $p.Find.Execute("lla tavalla.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lla tavalla (Signal jne.).", 2)

# ---------------------------------------------------------------------------
# Hunk 3: rewrite the IP field sentence and drop the _GoBack bookmark here
# (it will be re-created after this point once the new text exists).
# ---------------------------------------------------------------------------
$r = $d.Content
$old = " kenttään jossa sitä kysytään(tämä ei ole pakollista jos haluat toimittaa salatun viestitiedoston mulla sovelluksilla) "
$new = " kenttään jossa sitä kysytään ja vastaanottajan tietoliikenneportti. "
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, `
    $false, $new, 2)

# ---------------------------------------------------------------------------
# Hunk 4: rewrite the "Salattu viesti messageFile.txt ..." paragraph.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(12).Range
$old = "Salattu viesti messageFile.txt kopiodaan myös alussa valitsemaasi kansioon missä ovat muut salausavaimet. Jos haluat toimittaa viestin muulla tavalla kuin sovelluksella voit kopioda tiedoston tästä kansiosta.  Jos lähetät viestin sovelluksella tuhoa messageFile.txt. Jos taas valitsen oman tavan toimittaa messageFile.txt tiedoston on tiedosto syytä tuhota manuaalisesti siitä kansiosta johon se kopiotiin sen jälkeen kun olet lähettänyt tiedoston."
$new = "Salattu viesti messageFile.txt kopioidaan alussa valitsemaasi kansioon missä ovat muut salausavaimet. Jos haluat toimittaa viestin muulla tavalla kuin sovelluksella voit kopioida tiedoston tästä kansiosta tarvittaessa ja sulkea Viestittely sovelluksen ennen kuin viesti on lähetetty.  Jos lähetät viestin sovelluksella tuhoa messageFile.txt kun sovellus on saanut lähetyksen valmiiksi.  Jos taas valitsen oman tavan toimittaa messageFile.txt tiedoston on tiedosto syytä tuhota manuaalisesti siitä kansiosta johon se kopiotiin sen jälkeen kun olet lähettänyt tiedoston muulla tavalla (Signal jne.)."
$found = $p.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, `
    $false, $new, 2)

# ---------------------------------------------------------------------------
# Hunk 5 (implicit): moving the _GoBack bookmark below also removes it from
# its old spot inside "Viestin vastaanotto. Siirry ...". Re-adding a
# bookmark with the same name relocates it (Word bookmark names are unique).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("tietoliikenneportti. ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# Hunk 6 + 7: split the last content paragraph, add the new intro sentences
# about listening on the IP port, and append the closing sentences about
# reading the message once.
# ---------------------------------------------------------------------------
$r = $d.Content
$old = "Viestin vastaanoton aluksi valitse selaa nappi ja valitse kansio missä ovat"
$new = "^pKun viesti vastaanotetaan on asetettava internet-yhteyden kuuntelu aktiiviseksi IP:n kuuntelu painikkeella.   Viestin saavuttua aluksi valitse selaa nappi ja valitse kansio missä ovat"
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, `
    $false, $new, 2)

$r = $d.Content
$old = " viestikenttään."
$new = " viestikenttään. Viestin voi lukea vain kerran jonka jölkeen se tuhotaan. Jos haluat tallettaa avatun viestin voit kopioida sen viestikentästä."
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, `
    $false, $new, 2)
